# Applies the "Report notes.docx" edit described by the commit:
#   "Collapsed ASA categories, started changing table to only get
#    expected death rate for pd 39"
#
# Summary of structural changes to word/document.xml:
#   1. Paragraph 3 ("Data cleaning: ...removed proced = 2") gets a new
#      trailing run appended: ", categorized ASA into 1/2/3 and 4/5 b/c
#      no deaths in ASA=1 group"
#   2. The old paragraph 6 ("Modeled logistic regressions for death30 ~
#      ...") is replaced with new text ("Note: hospital 30 does not have
#      expected death rates ... pd 39, so its probabilities couldn't be
#      predicted") and keeps the _GoBack bookmark at its end.
#   3. The original "Modeled logistic regressions ..." text is kept, but
#      moved into a brand-new paragraph right after the Note paragraph
#      (without the bookmark).
#   4. A new trailing paragraph "Bootstrapping:  " + a tab character is
#      added right after the "Complete cases ..." paragraph.

$d = $word.ActiveDocument
$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$curlyApos = [char]0x2019

# ---------------------------------------------------------------------
# 1. "Data cleaning..." paragraph: append new run with the ASA note.
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$p3EndRange = $d.Range($p3.Range.End - 1, $p3.Range.End - 1)
$p3EndRange.InsertAfter(", categorized ASA into 1/2/3 and 4/5 b/c no deaths in ASA=1 group")

# ---------------------------------------------------------------------
# 2. Remove the whole old "Modeled logistic regressions..." paragraph;
#    we'll reinsert its text (unchanged) as its own new paragraph below,
#    and put new replacement text where it used to be.
# ---------------------------------------------------------------------
$oldModeledParagraphIndex = 6
$pOld = $d.Paragraphs($oldModeledParagraphIndex)
$pOld.Range.Delete()

# ---------------------------------------------------------------------
# 3. Insert the new "Note: hospital 30 ..." paragraph right after the
#    "Made table 1, table 2 ..." paragraph (index 5). Inserting a full
#    <w:p> fragment via InsertXML at a non-zero offset inside a
#    paragraph's range creates a brand-new paragraph immediately after
#    that paragraph without disturbing its own content.
# ---------------------------------------------------------------------
$pMadeTable = $d.Paragraphs(5)
$insertAfterMadeTable = $d.Range($pMadeTable.Range.Start + 5, $pMadeTable.Range.Start + 5)
$noteXml = '<w:p ' + $wordNs + '>' +
    '<w:r><w:t xml:space="preserve">Note: hospital 30 does not have expected death rates in table 2 b/c it had no BMI values for most recent </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>pd</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> 39, so its probabilities couldn' + $curlyApos + 't be predicted</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'
[void]$insertAfterMadeTable.InsertXML($noteXml)

# ---------------------------------------------------------------------
# 4. Insert the original "Modeled logistic regressions ..." text back in
#    as a fresh paragraph right after the new "Note:" paragraph.
# ---------------------------------------------------------------------
$pNote = $d.Paragraphs(6)
$insertAfterNote = $d.Range($pNote.Range.Start + 5, $pNote.Range.Start + 5)
$modeledXml = '<w:p ' + $wordNs + '>' +
    '<w:r><w:t xml:space="preserve">Modeled logistic regressions for death30 ~ </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>proced</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> + </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>asa</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> + </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>bmi</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> + albumin and then death30 ~ </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>proced</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> + </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>asa</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> + </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>bmi</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>; compared results (no different interpretations)</w:t></w:r>' +
    '</w:p>'
[void]$insertAfterNote.InsertXML($modeledXml)

# ---------------------------------------------------------------------
# 5. Add the new trailing "Bootstrapping:  " + tab paragraph right after
#    the "Complete cases ..." paragraph.
# ---------------------------------------------------------------------
$pComplete = $d.Paragraphs(8)
$insertAfterComplete = $d.Range($pComplete.Range.Start + 5, $pComplete.Range.Start + 5)
$bootstrapXml = '<w:p ' + $wordNs + '>' +
    '<w:r><w:t xml:space="preserve">Bootstrapping:  </w:t></w:r>' +
    '<w:r><w:tab/></w:r>' +
    '</w:p>'
[void]$insertAfterComplete.InsertXML($bootstrapXml)
